$d = $word.ActiveDocument

# 1) "...giving a short worked example demonstrating..." -> insert comma after "short"
$d.Content.Find.Execute(
    "giving a short worked example demonstrating",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "giving a short, worked example demonstrating", 2)

# 2) "...5 March 2020, and put this data on Vula..." -> remove comma
$d.Content.Find.Execute(
    "5 March 2020, and put this data on Vula",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "5 March 2020 and put this data on Vula", 2)

# 3) fix typo "assignement" -> "assignment"
$d.Content.Find.Execute(
    "assignement",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "assignment", 2)
